# Script produced by the betexplorer scraper run on 05-11-2023 14:45.
# 1) Two pairs of rows (21/22 and 51/52) had their match data (columns F:V)
#    swapped back into the correct order.
# 2) Two new matches (rows 81/82) were appended to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder the two mis-sequenced matches around rows 21/22 ---
$rangeA = $ws.Range("F21:V21")
$rangeB = $ws.Range("F22:V22")
$valsA = $rangeA.Value()
$valsB = $rangeB.Value()
$rangeA.Value = $valsB
$rangeB.Value = $valsA

# --- Reorder the two mis-sequenced matches around rows 51/52 ---
$rangeC = $ws.Range("F51:V51")
$rangeD = $ws.Range("F52:V52")
$valsC = $rangeC.Value()
$valsD = $rangeD.Value()
$rangeC.Value = $valsD
$rangeD.Value = $valsC

# --- Append the two newly scraped matches (rows 81 and 82) ---
# Copy the formatting (bold/centered index column, date-formatted E column)
# from the last existing data row so the new rows match the sheet's style.
$ws.Range("A80:V80").Copy() | Out-Null
$ws.Range("A81:V82").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A81").Value = 80
$ws.Range("B81").Value = "denmark"
$ws.Range("C81").Value = "superliga"
$ws.Range("D81").Value = "2023-2024"
$ws.Range("E81").Value = 45235.58333333334
$ws.Range("F81").Value = "Nordsjaelland"
$ws.Range("G81").Value = 1
$ws.Range("H81").Value = "Vejle"
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 1.33
$ws.Range("K81").Value = "30/10/2023 07:12"
$ws.Range("L81").Value = 1.3
$ws.Range("M81").Value = "05/11/2023 13:36"
$ws.Range("N81").Value = 5.51
$ws.Range("O81").Value = "30/10/2023 07:12"
$ws.Range("P81").Value = 5.81
$ws.Range("Q81").Value = "05/11/2023 13:55"
$ws.Range("R81").Value = 8.88
$ws.Range("S81").Value = "30/10/2023 07:12"
$ws.Range("T81").Value = 10.37
$ws.Range("U81").Value = "05/11/2023 13:55"
$ws.Range("V81").Value = "https://www.betexplorer.com/football/denmark/superliga/nordsjaelland-vejle/Sv8lDlMc/"

$ws.Range("A82").Value = 81
$ws.Range("B82").Value = "denmark"
$ws.Range("C82").Value = "superliga"
$ws.Range("D82").Value = "2023-2024"
$ws.Range("E82").Value = 45235.58333333334
$ws.Range("F82").Value = "Viborg"
$ws.Range("G82").Value = 2
$ws.Range("H82").Value = "Silkeborg"
$ws.Range("I82").Value = 1
$ws.Range("J82").Value = 2.45
$ws.Range("K82").Value = "29/10/2023 16:12"
$ws.Range("L82").Value = 2.57
$ws.Range("M82").Value = "05/11/2023 13:32"
$ws.Range("N82").Value = 3.45
$ws.Range("O82").Value = "29/10/2023 16:12"
$ws.Range("P82").Value = 3.55
$ws.Range("Q82").Value = "05/11/2023 13:09"
$ws.Range("R82").Value = 2.91
$ws.Range("S82").Value = "29/10/2023 16:12"
$ws.Range("T82").Value = 2.76
$ws.Range("U82").Value = "05/11/2023 13:32"
$ws.Range("V82").Value = "https://www.betexplorer.com/football/denmark/superliga/viborg-silkeborg/4pa98QzT/"
